# Auto-generated edit script: updates market-data derived cells (H:N)
# across the per-job profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed Universalis price snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 72726.5
$ws.Range("J3").Value = 72726.5
$ws.Range("L3").Value = 72726.5
$ws.Range("N3").Value = -72954.5

$ws.Range("H100").Value = 5853.7617
$ws.Range("I100").Value = 5853.7617
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5853.7617
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5312.7617
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 72726.5
$ws.Range("J102").Value = 72726.5
$ws.Range("L102").Value = 72726.5
$ws.Range("N102").Value = -79216.5

$ws.Range("H113").Value = 7592.737
$ws.Range("I113").Value = 7477.222
$ws.Range("J113").Value = 7696.7
$ws.Range("K113").Value = 7477.222
$ws.Range("L113").Value = 7696.7
$ws.Range("M113").Value = -4223.222
$ws.Range("N113").Value = -14204.7

$ws.Range("H116").Value = 20051.477
$ws.Range("I116").Value = 21776.295
$ws.Range("J116").Value = 12721
$ws.Range("K116").Value = 21776.295
$ws.Range("L116").Value = 12721
$ws.Range("M116").Value = -18334.295
$ws.Range("N116").Value = -19605

$ws.Range("H138").Value = 29769.73
$ws.Range("I138").Value = 2170.35
$ws.Range("K138").Value = 6511.049999999999
$ws.Range("M138").Value = -1371.049999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3673.074
$ws.Range("I2").Value = 3602.5833
$ws.Range("J2").Value = 4237
$ws.Range("K2").Value = 3602.5833
$ws.Range("L2").Value = 4237
$ws.Range("M2").Value = -3489.5833
$ws.Range("N2").Value = -4463

$ws.Range("H43").Value = 20894.273
$ws.Range("I43").Value = 22346.6
$ws.Range("K43").Value = 22346.6
$ws.Range("M43").Value = -22033.6

$ws.Range("H63").Value = 2841.7693
$ws.Range("J63").Value = 4749.5
$ws.Range("L63").Value = 4749.5
$ws.Range("N63").Value = -6121.5

$ws.Range("H66").Value = 2841.7693
$ws.Range("J66").Value = 4749.5
$ws.Range("L66").Value = 23747.5
$ws.Range("N66").Value = -30611.5

$ws.Range("H109").Value = 99888.5
$ws.Range("J109").Value = 99888.5
$ws.Range("L109").Value = 99888.5
$ws.Range("N109").Value = -102662.5

$ws.Range("H116").Value = 3673.074
$ws.Range("I116").Value = 3602.5833
$ws.Range("J116").Value = 4237
$ws.Range("K116").Value = 3602.5833
$ws.Range("L116").Value = 4237
$ws.Range("M116").Value = -1308.5833
$ws.Range("N116").Value = -8825

$ws.Range("H132").Value = 1068.8334
$ws.Range("I132").Value = 1068.8334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3206.5002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -676.5001999999999
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3673.074
$ws.Range("I3").Value = 3602.5833
$ws.Range("J3").Value = 4237
$ws.Range("K3").Value = 3602.5833
$ws.Range("L3").Value = 4237
$ws.Range("M3").Value = -3488.5833
$ws.Range("N3").Value = -4465

$ws.Range("H16").Value = 3169.3333
$ws.Range("I16").Value = 3169.3333
$ws.Range("K16").Value = 3169.3333
$ws.Range("M16").Value = -2999.3333

$ws.Range("H80").Value = 792.7059
$ws.Range("I80").Value = 1071.5714
$ws.Range("K80").Value = 1071.5714
$ws.Range("M80").Value = -73.57140000000004

$ws.Range("H83").Value = 792.7059
$ws.Range("I83").Value = 1071.5714
$ws.Range("K83").Value = 5357.857
$ws.Range("M83").Value = -365.857

$ws.Range("H86").Value = 1517.2222
$ws.Range("I86").Value = 1399
$ws.Range("J86").Value = 1551
$ws.Range("K86").Value = 1399
$ws.Range("L86").Value = 1551
$ws.Range("M86").Value = -276
$ws.Range("N86").Value = -3797

$ws.Range("H89").Value = 1517.2222
$ws.Range("I89").Value = 1399
$ws.Range("J89").Value = 1551
$ws.Range("K89").Value = 6995
$ws.Range("L89").Value = 7755
$ws.Range("M89").Value = -1379
$ws.Range("N89").Value = -18987

$ws.Range("H99").Value = 1172.75
$ws.Range("I99").Value = 1172.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1172.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 325.25
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 2181.3125
$ws.Range("I107").Value = 1793.6923
$ws.Range("J107").Value = 3861
$ws.Range("K107").Value = 1793.6923
$ws.Range("L107").Value = 3861
$ws.Range("M107").Value = 126.3077000000001
$ws.Range("N107").Value = -7701

$ws.Range("H134").Value = 2534.4443
$ws.Range("I134").Value = 2351.625
$ws.Range("K134").Value = 7054.875
$ws.Range("M134").Value = -4519.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4546965
$ws.Range("I31").Value = 6667583
$ws.Range("J31").Value = 2783.1428
$ws.Range("K31").Value = 6667583
$ws.Range("L31").Value = 2783.1428
$ws.Range("M31").Value = -6667288
$ws.Range("N31").Value = -3373.1428

$ws.Range("H34").Value = 4546965
$ws.Range("I34").Value = 6667583
$ws.Range("J34").Value = 2783.1428
$ws.Range("K34").Value = 6667583
$ws.Range("L34").Value = 2783.1428
$ws.Range("M34").Value = -6667381
$ws.Range("N34").Value = -3187.1428

$ws.Range("H107").Value = 1198.4166
$ws.Range("I107").Value = 1155.6666
$ws.Range("J107").Value = 1241.1666
$ws.Range("K107").Value = 1155.6666
$ws.Range("L107").Value = 1241.1666
$ws.Range("M107").Value = 764.3334
$ws.Range("N107").Value = -5081.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1900
$ws.Range("J117").Value = 1900
$ws.Range("L117").Value = 5700
$ws.Range("N117").Value = -12584

$ws.Range("H121").Value = 114513.586
$ws.Range("I121").Value = 152728.58
$ws.Range("J121").Value = 61012.6
$ws.Range("K121").Value = 458185.74
$ws.Range("L121").Value = 183037.8
$ws.Range("M121").Value = -456875.74
$ws.Range("N121").Value = -185657.8

$ws.Range("H122").Value = 1478.2
$ws.Range("I122").Value = 490
$ws.Range("K122").Value = 4410
$ws.Range("M122").Value = -1960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4637.2354
$ws.Range("J70").Value = 4620.125
$ws.Range("L70").Value = 4620.125
$ws.Range("N70").Value = -5160.125

$ws.Range("H73").Value = 4637.2354
$ws.Range("J73").Value = 4620.125
$ws.Range("L73").Value = 4620.125
$ws.Range("N73").Value = -6492.125

$ws.Range("H122").Value = 3448.361
$ws.Range("I122").Value = 3162.28
$ws.Range("K122").Value = 9486.84
$ws.Range("M122").Value = -7036.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4045.9375
$ws.Range("I7").Value = 4075.6667
$ws.Range("K7").Value = 4075.6667
$ws.Range("M7").Value = -3963.6667

$ws.Range("H126").Value = 4045.9375
$ws.Range("I126").Value = 4075.6667
$ws.Range("K126").Value = 12227.0001
$ws.Range("M126").Value = -9757.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4003
$ws.Range("J62").Value = 4003
$ws.Range("L62").Value = 4003
$ws.Range("N62").Value = -5251

$ws.Range("H65").Value = 4003
$ws.Range("J65").Value = 4003
$ws.Range("L65").Value = 20015
$ws.Range("N65").Value = -26255
